# startDeleteDialog.xlsx — apply the "cn181107" console upload edit.
#
# 1. Correct the capitalisation of the confirmation message in C2
#    ("delete" -> "Deletion", "instance termination" -> "Instance Termination").
#    Writing through .Value naturally retires the old shared-string entry and
#    appends the new text as a fresh shared string, which is what shifts the
#    "Cancel"/"OK" string indices the same way the source diff shows.
# 2. Widen column C so the longer English strings aren't clipped.
# 3. Move the selection to C16 (matches the saved cursor position in the diff).
# 4. Configure the print setup (portrait, paper size 9 / A4) for the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Confirm to enable Deletion on Instance Termination"

$ws.Columns("C").ColumnWidth = 40.33203125

$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

[void]$ws.Range("C16").Select()
